# TM52__TestJob1.xlsx - re-run log update
#
# The "readme" sheet keeps a running log (Table1, A1:E12) of every test
# run: index / Date / JobNo / sheet_name / Author. This edit represents a
# re-run of the test suite a day later (setup_class now records a fresh
# run), which:
#   - swaps the JobNo and Author columns (header + every data row), and
#   - bumps the logged run Date from 2022-02-23 to 2022-02-24.
# The "Project Information" sheet's "Date of Analysis" timestamp is
# likewise refreshed to match the new run.

$wb = $excel.ActiveWorkbook

# ---- "readme" sheet: Table1 log -------------------------------------
$readme = $wb.Worksheets.Item("readme")

# Header row: swap the JobNo / Author column titles (C1 <-> E1).
$readme.Cells.Item(1, 3).Value = "Author"
$readme.Cells.Item(1, 5).Value = "JobNo"

# Keep the table's own column-name metadata in sync with the header cells.
$table1 = $readme.ListObjects.Item("Table1")
$table1.ListColumns.Item(3).Name = "Author"
$table1.ListColumns.Item(5).Name = "JobNo"

# The new Date string ("20220224") reads as a plain number to Excel's
# usual Range.Value type coercion, which would silently store it as a
# numeric cell instead of text (unlike the original log, which is text).
# Route it through a throwaway text-formula cell + copy/PasteSpecial
# (values only) so it lands as a genuine string, then wipe the scratch
# cell so it leaves no trace in the sheet's used range.
$scratch = $readme.Cells.Item(20, 20)
$scratch.Formula = '="2022" & "0224"'
$scratch.Copy()

# Data rows 2-12: bump the Date, and swap the JobNo/Author values.
for ($row = 2; $row -le 12; $row++) {
    $readme.Cells.Item($row, 2).PasteSpecial(-4163)
    $readme.Cells.Item($row, 3).Value = "jovyan"
    $readme.Cells.Item($row, 5).Value = "/c/e"
}

$scratch.ClearContents()
$excel.CutCopyMode = $false

# ---- "Project Information" sheet: refresh the analysis timestamp ----
$projInfo = $wb.Worksheets.Item("Project Information")
$projInfo.Cells.Item(11, 2).Value = "2022-02-24 11:12:10.081416"
